# Append 6 new data rows (236-241) to Sheet1, mirroring format of row 235
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refRow = 235

# ---- Row 236 ----
# 1) Clone cell formatting (fill/border/font/number-format) from the row above
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(236,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(236).RowHeight = 16

# 2) Write the cell values. Text columns get NumberFormat "@" first so that
#    numeric-looking text (e.g. "2697", "2002") is stored as text, not a number.
#    (Columns B and C are intentionally left blank, matching the row above.)
$ws.Cells.Item(236,1).NumberFormat = "@"
$ws.Cells.Item(236,1).Value2 = '●'
$ws.Cells.Item(236,4).NumberFormat = "@"
$ws.Cells.Item(236,4).Value2 = '2697'
$ws.Cells.Item(236,5).NumberFormat = "@"
$ws.Cells.Item(236,5).Value2 = 'Event month'
$ws.Cells.Item(236,6).NumberFormat = "@"
$ws.Cells.Item(236,6).Value2 = '2: 2572'
$ws.Cells.Item(236,7).NumberFormat = "@"
$ws.Cells.Item(236,7).Value2 = '2: 2574'
$ws.Cells.Item(236,8).Value2 = 0
$ws.Cells.Item(236,9).NumberFormat = "@"
$ws.Cells.Item(236,9).Value2 = 'May'
$ws.Cells.Item(236,10).Value2 = 3
$ws.Cells.Item(236,11).Value2 = 0.011559
$ws.Cells.Item(236,12).NumberFormat = "@"
$ws.Cells.Item(236,12).Value2 = 'Sonia'
$ws.Cells.Item(236,13).NumberFormat = "@"
$ws.Cells.Item(236,13).Value2 = '11/14/18 11:31:00'

# 3) Re-apply the original formatting on top (restores the proper numFmt/style
#    index while keeping the text values just written).
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(236,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(236).RowHeight = 16

# 4) Formatting-only paste clears the values again, so re-apply A/D-M once more.
$ws.Cells.Item(236,1).NumberFormat = "@"
$ws.Cells.Item(236,1).Value2 = '●'
$ws.Cells.Item(236,4).NumberFormat = "@"
$ws.Cells.Item(236,4).Value2 = '2697'
$ws.Cells.Item(236,5).NumberFormat = "@"
$ws.Cells.Item(236,5).Value2 = 'Event month'
$ws.Cells.Item(236,6).NumberFormat = "@"
$ws.Cells.Item(236,6).Value2 = '2: 2572'
$ws.Cells.Item(236,7).NumberFormat = "@"
$ws.Cells.Item(236,7).Value2 = '2: 2574'
$ws.Cells.Item(236,8).Value2 = 0
$ws.Cells.Item(236,9).NumberFormat = "@"
$ws.Cells.Item(236,9).Value2 = 'May'
$ws.Cells.Item(236,10).Value2 = 3
$ws.Cells.Item(236,11).Value2 = 0.011559
$ws.Cells.Item(236,12).NumberFormat = "@"
$ws.Cells.Item(236,12).Value2 = 'Sonia'
$ws.Cells.Item(236,13).NumberFormat = "@"
$ws.Cells.Item(236,13).Value2 = '11/14/18 11:31:00'

# ---- Row 237 ----
# 1) Clone cell formatting (fill/border/font/number-format) from the row above
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(237,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(237).RowHeight = 16

# 2) Write the cell values. Text columns get NumberFormat "@" first so that
#    numeric-looking text (e.g. "2697", "2002") is stored as text, not a number.
#    (Columns B and C are intentionally left blank, matching the row above.)
$ws.Cells.Item(237,1).NumberFormat = "@"
$ws.Cells.Item(237,1).Value2 = '●'
$ws.Cells.Item(237,4).NumberFormat = "@"
$ws.Cells.Item(237,4).Value2 = '2697'
$ws.Cells.Item(237,5).NumberFormat = "@"
$ws.Cells.Item(237,5).Value2 = 'Event month'
$ws.Cells.Item(237,6).NumberFormat = "@"
$ws.Cells.Item(237,6).Value2 = '2: 2580'
$ws.Cells.Item(237,7).NumberFormat = "@"
$ws.Cells.Item(237,7).Value2 = '2: 2585'
$ws.Cells.Item(237,8).Value2 = 0
$ws.Cells.Item(237,9).NumberFormat = "@"
$ws.Cells.Item(237,9).Value2 = 'August'
$ws.Cells.Item(237,10).Value2 = 6
$ws.Cells.Item(237,11).Value2 = 0.023118
$ws.Cells.Item(237,12).NumberFormat = "@"
$ws.Cells.Item(237,12).Value2 = 'Sonia'
$ws.Cells.Item(237,13).NumberFormat = "@"
$ws.Cells.Item(237,13).Value2 = '11/14/18 11:31:00'

# 3) Re-apply the original formatting on top (restores the proper numFmt/style
#    index while keeping the text values just written).
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(237,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(237).RowHeight = 16

# 4) Formatting-only paste clears the values again, so re-apply A/D-M once more.
$ws.Cells.Item(237,1).NumberFormat = "@"
$ws.Cells.Item(237,1).Value2 = '●'
$ws.Cells.Item(237,4).NumberFormat = "@"
$ws.Cells.Item(237,4).Value2 = '2697'
$ws.Cells.Item(237,5).NumberFormat = "@"
$ws.Cells.Item(237,5).Value2 = 'Event month'
$ws.Cells.Item(237,6).NumberFormat = "@"
$ws.Cells.Item(237,6).Value2 = '2: 2580'
$ws.Cells.Item(237,7).NumberFormat = "@"
$ws.Cells.Item(237,7).Value2 = '2: 2585'
$ws.Cells.Item(237,8).Value2 = 0
$ws.Cells.Item(237,9).NumberFormat = "@"
$ws.Cells.Item(237,9).Value2 = 'August'
$ws.Cells.Item(237,10).Value2 = 6
$ws.Cells.Item(237,11).Value2 = 0.023118
$ws.Cells.Item(237,12).NumberFormat = "@"
$ws.Cells.Item(237,12).Value2 = 'Sonia'
$ws.Cells.Item(237,13).NumberFormat = "@"
$ws.Cells.Item(237,13).Value2 = '11/14/18 11:31:00'

# ---- Row 238 ----
# 1) Clone cell formatting (fill/border/font/number-format) from the row above
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(238,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(238).RowHeight = 16

# 2) Write the cell values. Text columns get NumberFormat "@" first so that
#    numeric-looking text (e.g. "2697", "2002") is stored as text, not a number.
#    (Columns B and C are intentionally left blank, matching the row above.)
$ws.Cells.Item(238,1).NumberFormat = "@"
$ws.Cells.Item(238,1).Value2 = '●'
$ws.Cells.Item(238,4).NumberFormat = "@"
$ws.Cells.Item(238,4).Value2 = '2697'
$ws.Cells.Item(238,5).NumberFormat = "@"
$ws.Cells.Item(238,5).Value2 = 'Event year'
$ws.Cells.Item(238,6).NumberFormat = "@"
$ws.Cells.Item(238,6).Value2 = '2: 2587'
$ws.Cells.Item(238,7).NumberFormat = "@"
$ws.Cells.Item(238,7).Value2 = '2: 2590'
$ws.Cells.Item(238,8).Value2 = 0
$ws.Cells.Item(238,9).NumberFormat = "@"
$ws.Cells.Item(238,9).Value2 = '2002'
$ws.Cells.Item(238,10).Value2 = 4
$ws.Cells.Item(238,11).Value2 = 0.015412
$ws.Cells.Item(238,12).NumberFormat = "@"
$ws.Cells.Item(238,12).Value2 = 'Sonia'
$ws.Cells.Item(238,13).NumberFormat = "@"
$ws.Cells.Item(238,13).Value2 = '11/14/18 11:31:00'

# 3) Re-apply the original formatting on top (restores the proper numFmt/style
#    index while keeping the text values just written).
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(238,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(238).RowHeight = 16

# 4) Formatting-only paste clears the values again, so re-apply A/D-M once more.
$ws.Cells.Item(238,1).NumberFormat = "@"
$ws.Cells.Item(238,1).Value2 = '●'
$ws.Cells.Item(238,4).NumberFormat = "@"
$ws.Cells.Item(238,4).Value2 = '2697'
$ws.Cells.Item(238,5).NumberFormat = "@"
$ws.Cells.Item(238,5).Value2 = 'Event year'
$ws.Cells.Item(238,6).NumberFormat = "@"
$ws.Cells.Item(238,6).Value2 = '2: 2587'
$ws.Cells.Item(238,7).NumberFormat = "@"
$ws.Cells.Item(238,7).Value2 = '2: 2590'
$ws.Cells.Item(238,8).Value2 = 0
$ws.Cells.Item(238,9).NumberFormat = "@"
$ws.Cells.Item(238,9).Value2 = '2002'
$ws.Cells.Item(238,10).Value2 = 4
$ws.Cells.Item(238,11).Value2 = 0.015412
$ws.Cells.Item(238,12).NumberFormat = "@"
$ws.Cells.Item(238,12).Value2 = 'Sonia'
$ws.Cells.Item(238,13).NumberFormat = "@"
$ws.Cells.Item(238,13).Value2 = '11/14/18 11:31:00'

# ---- Row 239 ----
# 1) Clone cell formatting (fill/border/font/number-format) from the row above
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(239,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(239).RowHeight = 16

# 2) Write the cell values. Text columns get NumberFormat "@" first so that
#    numeric-looking text (e.g. "2697", "2002") is stored as text, not a number.
#    (Columns B and C are intentionally left blank, matching the row above.)
$ws.Cells.Item(239,1).NumberFormat = "@"
$ws.Cells.Item(239,1).Value2 = '●'
$ws.Cells.Item(239,4).NumberFormat = "@"
$ws.Cells.Item(239,4).Value2 = '3651'
$ws.Cells.Item(239,5).NumberFormat = "@"
$ws.Cells.Item(239,5).Value2 = 'Event year'
$ws.Cells.Item(239,6).NumberFormat = "@"
$ws.Cells.Item(239,6).Value2 = '6: 665'
$ws.Cells.Item(239,7).NumberFormat = "@"
$ws.Cells.Item(239,7).Value2 = '6: 668'
$ws.Cells.Item(239,8).Value2 = 0
$ws.Cells.Item(239,9).NumberFormat = "@"
$ws.Cells.Item(239,9).Value2 = '2013'
$ws.Cells.Item(239,10).Value2 = 4
$ws.Cells.Item(239,11).Value2 = 0.008877
$ws.Cells.Item(239,12).NumberFormat = "@"
$ws.Cells.Item(239,12).Value2 = 'Sonia'
$ws.Cells.Item(239,13).NumberFormat = "@"
$ws.Cells.Item(239,13).Value2 = '11/14/18 11:33:00'

# 3) Re-apply the original formatting on top (restores the proper numFmt/style
#    index while keeping the text values just written).
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(239,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(239).RowHeight = 16

# 4) Formatting-only paste clears the values again, so re-apply A/D-M once more.
$ws.Cells.Item(239,1).NumberFormat = "@"
$ws.Cells.Item(239,1).Value2 = '●'
$ws.Cells.Item(239,4).NumberFormat = "@"
$ws.Cells.Item(239,4).Value2 = '3651'
$ws.Cells.Item(239,5).NumberFormat = "@"
$ws.Cells.Item(239,5).Value2 = 'Event year'
$ws.Cells.Item(239,6).NumberFormat = "@"
$ws.Cells.Item(239,6).Value2 = '6: 665'
$ws.Cells.Item(239,7).NumberFormat = "@"
$ws.Cells.Item(239,7).Value2 = '6: 668'
$ws.Cells.Item(239,8).Value2 = 0
$ws.Cells.Item(239,9).NumberFormat = "@"
$ws.Cells.Item(239,9).Value2 = '2013'
$ws.Cells.Item(239,10).Value2 = 4
$ws.Cells.Item(239,11).Value2 = 0.008877
$ws.Cells.Item(239,12).NumberFormat = "@"
$ws.Cells.Item(239,12).Value2 = 'Sonia'
$ws.Cells.Item(239,13).NumberFormat = "@"
$ws.Cells.Item(239,13).Value2 = '11/14/18 11:33:00'

# ---- Row 240 ----
# 1) Clone cell formatting (fill/border/font/number-format) from the row above
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(240,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(240).RowHeight = 16

# 2) Write the cell values. Text columns get NumberFormat "@" first so that
#    numeric-looking text (e.g. "2697", "2002") is stored as text, not a number.
#    (Columns B and C are intentionally left blank, matching the row above.)
$ws.Cells.Item(240,1).NumberFormat = "@"
$ws.Cells.Item(240,1).Value2 = '●'
$ws.Cells.Item(240,4).NumberFormat = "@"
$ws.Cells.Item(240,4).Value2 = '3651'
$ws.Cells.Item(240,5).NumberFormat = "@"
$ws.Cells.Item(240,5).Value2 = 'Event year'
$ws.Cells.Item(240,6).NumberFormat = "@"
$ws.Cells.Item(240,6).Value2 = '6: 655'
$ws.Cells.Item(240,7).NumberFormat = "@"
$ws.Cells.Item(240,7).Value2 = '6: 658'
$ws.Cells.Item(240,8).Value2 = 0
$ws.Cells.Item(240,9).NumberFormat = "@"
$ws.Cells.Item(240,9).Value2 = '2012'
$ws.Cells.Item(240,10).Value2 = 4
$ws.Cells.Item(240,11).Value2 = 0.008877
$ws.Cells.Item(240,12).NumberFormat = "@"
$ws.Cells.Item(240,12).Value2 = 'Sonia'
$ws.Cells.Item(240,13).NumberFormat = "@"
$ws.Cells.Item(240,13).Value2 = '11/14/18 11:33:00'

# 3) Re-apply the original formatting on top (restores the proper numFmt/style
#    index while keeping the text values just written).
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(240,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(240).RowHeight = 16

# 4) Formatting-only paste clears the values again, so re-apply A/D-M once more.
$ws.Cells.Item(240,1).NumberFormat = "@"
$ws.Cells.Item(240,1).Value2 = '●'
$ws.Cells.Item(240,4).NumberFormat = "@"
$ws.Cells.Item(240,4).Value2 = '3651'
$ws.Cells.Item(240,5).NumberFormat = "@"
$ws.Cells.Item(240,5).Value2 = 'Event year'
$ws.Cells.Item(240,6).NumberFormat = "@"
$ws.Cells.Item(240,6).Value2 = '6: 655'
$ws.Cells.Item(240,7).NumberFormat = "@"
$ws.Cells.Item(240,7).Value2 = '6: 658'
$ws.Cells.Item(240,8).Value2 = 0
$ws.Cells.Item(240,9).NumberFormat = "@"
$ws.Cells.Item(240,9).Value2 = '2012'
$ws.Cells.Item(240,10).Value2 = 4
$ws.Cells.Item(240,11).Value2 = 0.008877
$ws.Cells.Item(240,12).NumberFormat = "@"
$ws.Cells.Item(240,12).Value2 = 'Sonia'
$ws.Cells.Item(240,13).NumberFormat = "@"
$ws.Cells.Item(240,13).Value2 = '11/14/18 11:33:00'

# ---- Row 241 ----
# 1) Clone cell formatting (fill/border/font/number-format) from the row above
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(241,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(241).RowHeight = 16

# 2) Write the cell values. Text columns get NumberFormat "@" first so that
#    numeric-looking text (e.g. "2697", "2002") is stored as text, not a number.
#    (Columns B and C are intentionally left blank, matching the row above.)
$ws.Cells.Item(241,1).NumberFormat = "@"
$ws.Cells.Item(241,1).Value2 = '●'
$ws.Cells.Item(241,4).NumberFormat = "@"
$ws.Cells.Item(241,4).Value2 = '3910'
$ws.Cells.Item(241,5).NumberFormat = "@"
$ws.Cells.Item(241,5).Value2 = 'Event year'
$ws.Cells.Item(241,6).NumberFormat = "@"
$ws.Cells.Item(241,6).Value2 = '4: 2261'
$ws.Cells.Item(241,7).NumberFormat = "@"
$ws.Cells.Item(241,7).Value2 = '4: 2264'
$ws.Cells.Item(241,8).Value2 = 0
$ws.Cells.Item(241,9).NumberFormat = "@"
$ws.Cells.Item(241,9).Value2 = '2006'
$ws.Cells.Item(241,10).Value2 = 4
$ws.Cells.Item(241,11).Value2 = 0.020517
$ws.Cells.Item(241,12).NumberFormat = "@"
$ws.Cells.Item(241,12).Value2 = 'Sonia'
$ws.Cells.Item(241,13).NumberFormat = "@"
$ws.Cells.Item(241,13).Value2 = '11/14/18 11:35:00'

# 3) Re-apply the original formatting on top (restores the proper numFmt/style
#    index while keeping the text values just written).
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item($refRow,$c).Copy()
    $ws.Cells.Item(241,$c).PasteSpecial(-4122)
}
$ws.Rows.Item(241).RowHeight = 16

# 4) Formatting-only paste clears the values again, so re-apply A/D-M once more.
$ws.Cells.Item(241,1).NumberFormat = "@"
$ws.Cells.Item(241,1).Value2 = '●'
$ws.Cells.Item(241,4).NumberFormat = "@"
$ws.Cells.Item(241,4).Value2 = '3910'
$ws.Cells.Item(241,5).NumberFormat = "@"
$ws.Cells.Item(241,5).Value2 = 'Event year'
$ws.Cells.Item(241,6).NumberFormat = "@"
$ws.Cells.Item(241,6).Value2 = '4: 2261'
$ws.Cells.Item(241,7).NumberFormat = "@"
$ws.Cells.Item(241,7).Value2 = '4: 2264'
$ws.Cells.Item(241,8).Value2 = 0
$ws.Cells.Item(241,9).NumberFormat = "@"
$ws.Cells.Item(241,9).Value2 = '2006'
$ws.Cells.Item(241,10).Value2 = 4
$ws.Cells.Item(241,11).Value2 = 0.020517
$ws.Cells.Item(241,12).NumberFormat = "@"
$ws.Cells.Item(241,12).Value2 = 'Sonia'
$ws.Cells.Item(241,13).NumberFormat = "@"
$ws.Cells.Item(241,13).Value2 = '11/14/18 11:35:00'

$excel.CutCopyMode = 0
Write-Host "Appended rows 236-241. UsedRange:" $ws.UsedRange.Address()
